{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Texts of the two footer-ish paragraphs that must be removed, expected to\n// immediately follow the bibliography's closing entry (with one blank\n// paragraph separating them from it).\nconst verNoJupiter = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightLine =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Find the bibliography entry paragraph ending in \"...2016;\" \u2014 the anchor\n// right before the three paragraphs that need to go away.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Educations (3\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (anchorIndex !== -1) {\n  const blank = items[anchorIndex + 1];\n  const verPara = items[anchorIndex + 2];\n  const copyrightPara = items[anchorIndex + 3];\n\n  if (blank && blank.text === \"\") {\n    toDelete.push(blank);\n  }\n  if (verPara && verPara.text === verNoJupiter) {\n    toDelete.push(verPara);\n  }\n  if (copyrightPara && copyrightPara.text === copyrightLine) {\n    toDelete.push(copyrightPara);\n  }\n}\n\n// Delete from the end backwards so earlier indices/objects stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the bibliography entry paragraph that ends the reference list:\n# \"...Educations (3\u00aa Edi\u00e7\u00e3o), 872 p., 2016;\"\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Educations (3*Edi*2016*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    # Immediately after that paragraph sit three paragraphs that must go:\n    #   1) a blank paragraph\n    #   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n    #   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n    #      pages. Original theme under Creative Commons Attribution\"\n    # Deleting the paragraph right after the anchor, three times in a row,\n    # removes exactly those three (each delete shifts the following\n    # paragraphs up into that slot) while leaving the anchor paragraph and\n    # the trailing blank / page-break paragraphs untouched.\n    for ($n = 0; $n -lt 3; $n++) {\n        $d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n    }\n}\n"}
